$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each touched cell, force text format before assigning so that
# numeric-looking values (e.g. "0.9993") are written as Text, matching
# the original inline-string cell type, instead of being auto-coerced
# into Number cells by COM type inference.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.089.53'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.849.63'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.6925'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -5.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '237.85'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9995'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07730'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +8.65%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3033'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.23'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08111'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.860.98'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.37%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.200'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.02'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.091.27'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.736'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -4.44%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007770'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.18'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '235.42'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.0000'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.097.30'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9996'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.588'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.971'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.44%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.39'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1426'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -7.33%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.977'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.397'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.493'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.486'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.31%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.015'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05219'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.41%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.180'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.09%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7009'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.34%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.022'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.653'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.675'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9139'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +5.51%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.089.89'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.003'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.42%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '70.61'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.78%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9992'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.84'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.766'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.58%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.994.71'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.143'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.960'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -6.32%  '
